$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (bold) ----
$ws.Range("A1").Value = "Progress"
$ws.Range("B1").Value = "Assignee"
$ws.Range("C1").Value = "User Role"
$ws.Range("D1").Value = "Feature"
$ws.Range("E1").Value = "ETA"
$ws.Range("F1").Value = "DEADLINE"
$ws.Range("A1:F1").Font.Bold = $true

# ---- Data rows ----
$ws.Range("A2").Value = "In Progress"
$ws.Range("B2").Value = "Carson"
$ws.Range("C2").Value = "Manager"
$ws.Range("D2").Value = "Turn records of items tracked by employees into a table on the desktop app "
$ws.Range("E2").Value = 43900
$ws.Range("F2").Value = 43901

$ws.Range("A3").Value = "In Progress"
$ws.Range("B3").Value = "Tristen"
$ws.Range("C3").Value = "Member"
$ws.Range("D3").Value = "Be able to remove and update addresses "
$ws.Range("E3").Value = 43900
$ws.Range("F3").Value = 43901

$ws.Range("A4").Value = "In Progress"
$ws.Range("B4").Value = "Luke"
$ws.Range("C4").Value = "Warehouse employee "
$ws.Range("D4").Value = "Add ability to flag for damaged returns "
$ws.Range("E4").Value = 43900
$ws.Range("F4").Value = 43901

$ws.Range("A5").Value = "In Progress"
$ws.Range("B5").Value = "Luke/Tristen/Carson"
$ws.Range("C5").Value = "Librarian"
$ws.Range("D5").Value = "Implement view for librarians  "
$ws.Range("E5").Value = 43900
$ws.Range("F5").Value = 43901

$ws.Range("A6").Value = "In Progress"
$ws.Range("B6").Value = "Luke/Tristen/Carson"
$ws.Range("C6").Value = "Librarian"
$ws.Range("D6").Value = "implement view for viewing users history/information "
$ws.Range("E6").Value = 43900
$ws.Range("F6").Value = 43901

# ---- Formatting ----
# Feature column uses a smaller Arial font
$ws.Range("D2:D6").Font.Name = "Arial"
$ws.Range("D2:D6").Font.Size = 10

# ETA / DEADLINE columns formatted as dates
$ws.Range("E2:F6").NumberFormat = "d-mmm"

# ---- View state ----
$excel.ActiveWindow.Zoom = 145
$ws.Range("D7").Select()
